$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Z10").Value = "x"
$ws.Range("Z10").Font.FontStyle = "Bold"
Write-Host "done"
